# Update the multiplication-table answers in the single table of the
# document. The table has 5 "answer" rows (table rows 1, 5, 10, 15, 20,
# each with 5 cells) separated by blank rows; every answer cell's text
# is replaced with its new value, preserving the existing run/paragraph
# formatting (font, size, alignment) already applied to each cell.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row 1 (was: 32x11=352, 99x26=2574, 63x14=882, 17x22=374, 19x86=1634)
$t.Cell(1,1).Range.Text = "56×88=4928"
$t.Cell(1,2).Range.Text = "35×13=455"
$t.Cell(1,3).Range.Text = "16×28=448"
$t.Cell(1,4).Range.Text = "54×92=4968"
$t.Cell(1,5).Range.Text = "89×16=1424"

# Row 5 (was: 91x69=6279, 56x77=4312, 14x17=238, 67x51=3417, 49x39=1911)
$t.Cell(5,1).Range.Text = "60×98=5880"
$t.Cell(5,2).Range.Text = "21×49=1029"
$t.Cell(5,3).Range.Text = "99×65=6435"
$t.Cell(5,4).Range.Text = "27×40=1080"
$t.Cell(5,5).Range.Text = "99×19=1881"

# Row 10 (was: 55x76=4180, 72x72=5184, 16x19=304, 93x74=6882, 30x59=1770)
$t.Cell(10,1).Range.Text = "17×79=1343"
$t.Cell(10,2).Range.Text = "89×11=979"
$t.Cell(10,3).Range.Text = "28×58=1624"
$t.Cell(10,4).Range.Text = "30×57=1710"
$t.Cell(10,5).Range.Text = "70×28=1960"

# Row 15 (was: 24x85=2040, 44x92=4048, 16x82=1312, 63x74=4662, 30x19=570)
$t.Cell(15,1).Range.Text = "17×28=476"
$t.Cell(15,2).Range.Text = "94×39=3666"
$t.Cell(15,3).Range.Text = "50×38=1900"
$t.Cell(15,4).Range.Text = "51×64=3264"
$t.Cell(15,5).Range.Text = "46×69=3174"

# Row 20 (was: 81x37=2997, 63x74=4662, 65x41=2665, 73x31=2263, 60x40=2400)
$t.Cell(20,1).Range.Text = "68×92=6256"
$t.Cell(20,2).Range.Text = "29×21=609"
$t.Cell(20,3).Range.Text = "98×26=2548"
$t.Cell(20,4).Range.Text = "23×78=1794"
$t.Cell(20,5).Range.Text = "85×77=6545"
